$d = $word.ActiveDocument
$full = $d.Range(286, 351)
Write-Host "Before:[$($full.Text)]"
$full.Text = "#CEP_CLIENTE na cidade de #CIDADE_CLIENTE - #SIGLA_ESTADO_CLIENTE"
Write-Host "After:[$($full.Text)] Start:" $full.Start "End:" $full.End

# boundaries relative to start=286: 12, 26, 41, 44 (end=65)
$b1 = $d.Range(286+12, 286+65)  # " na cidade de #CIDADE_CLIENTE - #SIGLA_ESTADO_CLIENTE"
Write-Host "b1 before toggling Bold:" $b1.Font.Bold "Text:[$($b1.Text)]"
$b1.Font.Bold = 1
$b1.Font.Bold = 0

$b2 = $d.Range(286+26, 286+65)  # "#CIDADE_CLIENTE - #SIGLA_ESTADO_CLIENTE"
$b2.Font.Bold = 1
$b2.Font.Bold = 0

$b3 = $d.Range(286+41, 286+65)  # " - #SIGLA_ESTADO_CLIENTE"
$b3.Font.Bold = 1
$b3.Font.Bold = 0

$b4 = $d.Range(286+44, 286+65)  # "#SIGLA_ESTADO_CLIENTE"
$b4.Font.Bold = 1
$b4.Font.Bold = 0

Write-Host "done"
